$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 249.28572
$ws.Range("I38").Value = 149
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 447
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = -75
$ws.Range("N38").Value = -2244

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1478.3334
$ws.Range("I58").Value = 967.5
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 2902.5
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -2752.5
$ws.Range("N58").Value = -7800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1749
$ws.Range("I94").Value = 1749
$ws.Range("K94").Value = 1749
$ws.Range("M94").Value = -1298

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2984.7144
$ws.Range("I138").Value = 1158
$ws.Range("J138").Value = 3999.5557
$ws.Range("K138").Value = 3474
$ws.Range("L138").Value = 11998.6671
$ws.Range("M138").Value = 1666
$ws.Range("N138").Value = -22278.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4925.1113
$ws.Range("I32").Value = 3481.7273
$ws.Range("J32").Value = 11276
$ws.Range("K32").Value = 3481.7273
$ws.Range("L32").Value = 11276
$ws.Range("M32").Value = -3194.7273
$ws.Range("N32").Value = -11850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4107.077
$ws.Range("I110").Value = 3115.6667
$ws.Range("J110").Value = 4956.857
$ws.Range("K110").Value = 3115.6667
$ws.Range("L110").Value = 4956.857
$ws.Range("M110").Value = -1070.6667
$ws.Range("N110").Value = -9046.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1922.3077
$ws.Range("I132").Value = 1922.3077
$ws.Range("K132").Value = 5766.9231
$ws.Range("M132").Value = -3236.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4429.75
$ws.Range("I134").Value = 2812.6667
$ws.Range("K134").Value = 8438.000100000001
$ws.Range("M134").Value = -5903.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 225.08333
$ws.Range("I7").Value = 247.42857
$ws.Range("J7").Value = 193.8
$ws.Range("K7").Value = 247.42857
$ws.Range("L7").Value = 193.8
$ws.Range("M7").Value = -134.42857
$ws.Range("N7").Value = -419.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20962.268
$ws.Range("I31").Value = 12716.333
$ws.Range("J31").Value = 33331.168
$ws.Range("K31").Value = 12716.333
$ws.Range("L31").Value = 33331.168
$ws.Range("M31").Value = -12421.333
$ws.Range("N31").Value = -33921.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 20962.268
$ws.Range("I34").Value = 12716.333
$ws.Range("J34").Value = 33331.168
$ws.Range("K34").Value = 12716.333
$ws.Range("L34").Value = 33331.168
$ws.Range("M34").Value = -12514.333
$ws.Range("N34").Value = -33735.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4150
$ws.Range("I105").Value = 533.3333
$ws.Range("K105").Value = 533.3333
$ws.Range("M105").Value = 1213.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 267.81818
$ws.Range("I107").Value = 143
$ws.Range("K107").Value = 143
$ws.Range("M107").Value = 1777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1759.7142
$ws.Range("I122").Value = 1429.5
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 4288.5
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -1838.5
$ws.Range("N122").Value = -11500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7201.6665
$ws.Range("I132").Value = 7201.6665
$ws.Range("K132").Value = 21604.9995
$ws.Range("M132").Value = -19074.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 55564.61
$ws.Range("I4").Value = 9.538462000000001
$ws.Range("K4").Value = 28.615386
$ws.Range("M4").Value = 83.384614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 5301.875
$ws.Range("J94").Value = 5230.7144
$ws.Range("L94").Value = 15692.1432
$ws.Range("N94").Value = -17044.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 428.5
$ws.Range("I109").Value = 428.5
$ws.Range("K109").Value = 1285.5
$ws.Range("M109").Value = -245.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5959.7144
$ws.Range("I132").Value = 5619.8335
$ws.Range("K132").Value = 16859.5005
$ws.Range("M132").Value = -14329.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5053
$ws.Range("I7").Value = 2856
$ws.Range("K7").Value = 2856
$ws.Range("M7").Value = -2744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5230.2
$ws.Range("I40").Value = 5230.2
$ws.Range("K40").Value = 5230.2
$ws.Range("M40").Value = -5094.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6799.8
$ws.Range("I61").Value = 6799.8
$ws.Range("K61").Value = 6799.8
$ws.Range("M61").Value = -6597.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1430.3636
$ws.Range("I93").Value = 1430.3636
$ws.Range("K93").Value = 1430.3636
$ws.Range("M93").Value = -182.3635999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6799.8
$ws.Range("I113").Value = 6799.8
$ws.Range("K113").Value = 6799.8
$ws.Range("M113").Value = -4629.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6941.75
$ws.Range("I122").Value = 6941.75
$ws.Range("K122").Value = 20825.25
$ws.Range("M122").Value = -18375.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5053
$ws.Range("I126").Value = 2856
$ws.Range("K126").Value = 8568
$ws.Range("M126").Value = -6098

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 27308.334
$ws.Range("I132").Value = 27253.572
$ws.Range("J132").Value = 27500
$ws.Range("K132").Value = 81760.716
$ws.Range("L132").Value = 82500
$ws.Range("M132").Value = -79230.716
$ws.Range("N132").Value = -87560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3250
$ws.Range("J62").Value = 3250
$ws.Range("L62").Value = 3250
$ws.Range("N62").Value = -4498

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3250
$ws.Range("J65").Value = 3250
$ws.Range("L65").Value = 16250
$ws.Range("N65").Value = -22490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1914.4286
$ws.Range("I107").Value = 1100.5
$ws.Range("J107").Value = 2999.6667
$ws.Range("K107").Value = 3301.5
$ws.Range("L107").Value = 8999.000100000001
$ws.Range("M107").Value = -1381.5
$ws.Range("N107").Value = -12839.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3780.8
$ws.Range("I132").Value = 3780.8
$ws.Range("K132").Value = 11342.4
$ws.Range("M132").Value = -8812.400000000001
